$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing commas from tribunal court address fields
$ws.Range("B3").Value = "Manchester Employment Tribunal"
$ws.Range("B4").Value = "Alexandra House"
$ws.Range("B5").Value = "14-22 The Parsonage"
$ws.Range("B6").Value = "Manchester"
$ws.Range("B12").Value = "Eagle Building"
$ws.Range("B13").Value = "215 Bothwell Street"
$ws.Range("B14").Value = "Glasgow"

# Update the selected cell shown in the worksheet view
$ws.Range("B7").Select()
